$d = $word.ActiveDocument

# 1. "Daniel Garza, Norma Espinosa, " + "Meina"(spell) + " " + "Bian"(spell) + ", Adam Freedman, Angele Yazbec"
#    -> single run, no proofErr markers.
$d.Content.Find.Execute("Daniel Garza, Norma Espinosa, Meina Bian, Adam Freedman, Angele Yazbec", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Daniel Garza, Norma Espinosa, Meina Bian, Adam Freedman, Angele Yazbec", 2) | Out-Null

# 2. "Obtain coordi" + "nates for the countries" -> single run.
$d.Content.Find.Execute("Obtain coordinates for the countries", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Obtain coordinates for the countries", 2) | Out-Null

# 3. "Using HT" + "ML and JavaScript (" + "GeoMapping"(spell) + ", Leaflet, " + "Plotly"(spell) + ")" -> single run.
$d.Content.Find.Execute("Using HTML and JavaScript (GeoMapping, Leaflet, Plotly)", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Using HTML and JavaScript (GeoMapping, Leaflet, Plotly)", 2) | Out-Null

# 4. "Bar/scatter/line charts show the trends" + " over time" -> single run.
$d.Content.Find.Execute("Bar/scatter/line charts show the trends over time", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Bar/scatter/line charts show the trends over time", 2) | Out-Null

# 5. Actual text change: "Cleaning data in Pandas" -> "HTML"
$d.Content.Find.Execute("Cleaning data in Pandas", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "HTML", 2) | Out-Null

# 6. "Assist with visualizations " + "as needed" -> single run (applies harmlessly to the
#    already-merged later occurrences too, since find text == replace text).
$d.Content.Find.Execute("Assist with visualizations as needed", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Assist with visualizations as needed", 2) | Out-Null

# 7 & 8. Two paragraphs still carry <w:proofErr spellStart/> / <w:proofErr
#    spellEnd/> wrapping their only (or last) run, with no neighbouring run
#    on the far side of the tag for an ordinary Find/Replace to merge across:
#      - "HTML and " + "GeoMapping"(spell)         -> trailing spellEnd
#        sits right before </w:p> (no run follows it).
#      - "Meina"(spell) standalone paragraph        -> both spellStart and
#        spellEnd sit at the paragraph's edges (its only run has no
#        neighbour on either side).
#    Find/Replace only drops a proofErr tag when the matched span forces a
#    run-merge across the boundary where the tag sits. So for each such
#    paragraph we temporarily pad one space on whichever side(s) lack a
#    neighbouring run (creating one), run a single Find/Replace over the
#    padded text (which merges across every proofErr boundary in one go),
#    then strip the padding back off - leaving one clean run with no
#    proofErr wrapper left behind.
function Remove-ProofErrPadding($p, $core) {
    $start = $p.Range.Start
    $end = $p.Range.End - 1

    $lead = $d.Range($start, $start)
    $lead.InsertBefore(" ")

    $trail = $d.Range($end + 1, $end + 1)
    $trail.InsertAfter(" ")

    $padded = " " + $core + " "
    $p.Range.Find.Execute($padded, $false, $false, $false, $false, $false, $true, 1, $false, $padded, 2) | Out-Null
    $p.Range.Find.Execute($padded, $false, $false, $false, $false, $false, $true, 1, $false, $core, 2) | Out-Null
}

foreach ($p in $d.Paragraphs) {
    $ptext = $p.Range.Text
    if ($ptext -eq ("HTML and GeoMapping" + [char]13)) {
        Remove-ProofErrPadding $p "HTML and GeoMapping"
    } elseif ($ptext -eq ("Meina" + [char]13)) {
        Remove-ProofErrPadding $p "Meina"
    }
}
